{"js": "// Update the \"imagens(...)\" table definition:\n//   imagens(paciente_id,link,descricao)\n// becomes\n//   imagens(Id,paciente_id,arquivo,descricao)\n//\n// i.e. \"paciente_id\" -> \"Id\" (stays underlined),\n//      \"link\" -> \"paciente_id\" (stays underlined),\n//      a new \",arquivo\" field is inserted right before \",descricao\" (not\n//      underlined, matching the surrounding plain text),\n// and the document's \"_GoBack\" last-edit bookmark ends up sitting inside\n// the new \"paciente_id\" run (between \"pacient\" and \"e_id\"), matching\n// where Word would leave it after that edit.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"imagens(\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'imagens(...)' paragraph\");\n}\n\n// 1) paciente_id -> Id (first underlined field in this paragraph)\nlet hits = target.search(\"paciente_id\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nhits.items[0].insertText(\"Id\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) link -> paciente_id (still underlined)\nhits = target.search(\"link\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nhits.items[0].insertText(\"paciente_id\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) insert \",arquivo\" right before \",descricao\" (plain formatting,\n//    inherited from the existing \",descricao\" run, i.e. no underline)\nhits = target.search(\",descricao\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nhits.items[0].insertText(\",arquivo,descricao\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) relocate the \"_GoBack\" bookmark from wherever it currently sits to\n//    the middle of the \"paciente_id\" run we just created (right after\n//    \"pacient\", before \"e_id\") -- this is where Word leaves the\n//    last-edit-position bookmark after typing \"paciente_id\" over \"link\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nhits = target.search(\"pacient\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nconst afterPacient = hits.items[0].getRange(\"After\");\nafterPacient.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Update the \"imagens(...)\" table definition:\n#   imagens(paciente_id,link,descricao)\n# becomes\n#   imagens(Id,paciente_id,arquivo,descricao)\n#\n# i.e. \"paciente_id\" -> \"Id\" (stays underlined),\n#      \"link\" -> \"paciente_id\" (stays underlined),\n#      a new \",arquivo\" field is inserted right before \",descricao\" (plain\n#      formatting, no underline, like the rest of that text),\n# and the document's \"_GoBack\" last-edit bookmark ends up sitting inside\n# the new \"paciente_id\" run (between \"pacient\" and \"e_id\"), matching\n# where Word would leave it after that edit.\n\n$d = $word.ActiveDocument\n\n# Locate the \"imagens(...)\" paragraph.\n$paras = $d.Paragraphs\n$target = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($t.StartsWith(\"imagens(\")) {\n        $target = $paras.Item($i)\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the 'imagens(...)' paragraph\"\n}\n\n# 1) paciente_id -> Id (first underlined field in this paragraph)\n$rng = $target.Range\n$null = $rng.Find.Execute(\"paciente_id\", $false, $false, $false, $false, $false, $true, 1, $false, \"Id\", 2)\n\n# 2) link -> paciente_id (still underlined)\n$rng = $target.Range\n$null = $rng.Find.Execute(\"link\", $false, $false, $false, $false, $false, $true, 1, $false, \"paciente_id\", 2)\n\n# 3) insert \",arquivo\" right before \",descricao\" (plain formatting,\n#    inherited from the existing \",descricao\" text, i.e. no underline)\n$rng = $target.Range\n$null = $rng.Find.Execute(\",descricao\", $false, $false, $false, $false, $false, $true, 1, $false, \",arquivo,descricao\", 2)\n\n# 4) relocate the \"_GoBack\" bookmark from wherever it currently sits to\n#    the middle of the \"paciente_id\" run we just created (right after\n#    \"pacient\", before \"e_id\") -- this is where Word leaves the\n#    last-edit-position bookmark after typing \"paciente_id\" over \"link\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$findRng = $target.Range\n$null = $findRng.Find.Execute(\"pacient\")\n$bmRange = $d.Range($findRng.End, $findRng.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
